$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.573.19"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +8.50%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.582.52"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +10.19%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "504.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.20%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +9.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.626"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +24.91%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.581.60"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +10.34%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +15.93%  "

$ws.Range("E11").Value = "  +6.87%  "

$ws.Range("E12").Value = "  +6.96%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.127"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.97%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.030.57"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +10.18%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.381.32"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.45%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.80"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +9.13%  "

$ws.Range("E17").Value = "  +5.48%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.588.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +10.29%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "334.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.68%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.78%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +8.37%  "

$ws.Range("E23").Value = "  +0.94%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.64%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.414"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.86%  "

$ws.Range("E26").Value = "  +10.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.701.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.56%  "

$ws.Range("E28").Value = "  -0.29%  "

$ws.Range("E29").Value = "  +10.20%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.20%  "

$ws.Range("E31").Value = "  -0.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "158.06"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.54%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.32"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.49%  "

$ws.Range("E34").Value = "  +7.44%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.48"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.86%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.93"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +11.54%  "

$ws.Range("E37").Value = "  +9.82%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.851"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.54%  "

$ws.Range("E39").Value = "  +12.67%  "

$ws.Range("E40").Value = "  +9.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "35.13"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.40%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "289.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +15.63%  "

$ws.Range("E43").Value = "  +7.70%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.624"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +9.29%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0563"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.64%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.20%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.35"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +15.87%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0235"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.17%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.729"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +16.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.75"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.11%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.13%  "
